$d = $word.ActiveDocument

$search = $d.Content
$search.Find.Execute("interfacci agrafica") | Out-Null
$foundStart = $search.Start
$foundEnd = $search.End

# Step 1: set full replacement text first
$r = $d.Range($foundStart, $foundEnd)
$r.Text = "interfaccia grafica"

# Step 2: toggle bold on the space char then off
$spaceR = $d.Range($foundStart + 11, $foundStart + 12)
$spaceR.Bold = 1
$spaceR.Bold = 0

Write-Output ("final text=[" + $d.Range($foundStart, $foundStart+20).Text + "]")
